$wb = $excel.ActiveWorkbook

# --- Sheet1: update admin login cells ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("A2").Value = "Admin"

# --- Sheet2: add the "advanced topic" data row ---
$ws2 = $wb.Worksheets.Item("Sheet2")
# Write in this order so the shared-string table is built up in the same
# sequence as the target workbook (Dec, This is a test, Cassidy Hope, US - FMLA, ...)
$ws2.Range("C1").Value = "Dec"
$ws2.Range("I1").Value = "This is a test"
$ws2.Range("A1").Value = "Cassidy Hope"
$ws2.Range("B1").Value = "US - FMLA"
$ws2.Range("D1").Value = 2020
$ws2.Range("E1").Value = 14
$ws2.Range("F1").Value = "Dec"
$ws2.Range("G1").Value = 2020
$ws2.Range("H1").Value = 19

# Column widths (best-fit-ish), matching the target's custom widths as
# closely as this engine's character->pixel rounding allows.
$ws2.Columns.Item(1).ColumnWidth = 10.74
$ws2.Columns.Item(2).ColumnWidth = 17.9
$ws2.Columns.Item(7).ColumnWidth = 10.9

# Portrait page orientation for Sheet2.
$ws2.PageSetup.Orientation = 1

# --- Remove Sheet3 entirely ---
$wb.Worksheets.Item("Sheet3").Delete() | Out-Null

# --- Selections / active sheet ---
# Select A2 on Sheet1 first (leaves it as the non-active tab).
$ws1.Range("A2").Select() | Out-Null
# Then make Sheet2 the active tab, with H1 selected - matches the diff's
# final sheetView/selection + workbook activeTab state.
$ws2.Activate() | Out-Null
$ws2.Range("H1").Select() | Out-Null
